# "raaputin viimeistelty + tyotunnit lisatty"
# Finish off the scratch rows and add a new work-hours entry on Taul1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")
$ws.Activate()

# New work session on row 18: 15.5 (serial 41409), 6 tunnit, "projektihuoneella".
$ws.Range("A18").NumberFormat = "d-mmm"
$ws.Range("A18").Value = 41409
$ws.Range("B18").Value = 6
$ws.Range("C18").Value = "projektihuoneella"

# Total the logged hours into row 27.
$ws.Range("B27").Formula = "=SUM(B2:B26)"

# Leave the selection/scroll position where the author ended up editing.
$ws.Range("B28").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
